$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The exact order below reproduces the shared-string table ordering of the
# target workbook (new strings are interned in the order Excel first
# encounters them while the cells are written).

$ws.Range("A37").Value = "Virtual Security"
$ws.Range("C36").Value = "Runs inside of an operating system of a physical host machine. E.g. VirtualBox"
$ws.Range("C37").Value = "Security features, patch management, hardware maintenance, resilient and high availability"
$ws.Range("B37").Value = "Virtualization characteristics"
$ws.Range("B38").Value = "Network Separation"
$ws.Range("B39").Value = "Snapshots and backups"
$ws.Range("C38").Value = "Separate VMs on one server using virtual switch"
$ws.Range("C39").Value = "Easy to backup, save states, debug"
$ws.Range("B40").Value = "Virtual Threats"
$ws.Range("C40").Value = "Anything that can happen to VM can happen to physical machine"
$ws.Range("B41").Value = "VM Sprawl"
$ws.Range("C41").Value = "Various virtual machines with no centralization - bad thing"
$ws.Range("B42").Value = "VM Escape"
$ws.Range("C42").Value = "Escapes VM and affects host system"
$ws.Range("B43").Value = "Hardening Virtualization"
$ws.Range("C43").Value = "Remove remnant data, make good policies, define user privileges, patch everything"
$ws.Range("B44").Value = "Cloud Acess Security Broker"
$ws.Range("C44").Value = "Intermediary between local infrastructure and the cloud. Usually on the cloud. Watches for malware, and controls policies"
$ws.Range("A45").Value = "Containerisation"
$ws.Range("B45").Value = "Function"
$ws.Range("C45").Value = "Runs isolated instances of programs and services"
$ws.Range("C46").Value = "Self-contained application that communicates with network resources that are permitted"
$ws.Range("C47").Value = "Containers can depend on each other and can be configured to communicate on a single host"
$ws.Range("B46").Value = ">>"
$ws.Range("C48").Value = "Containers run a single program with all its dependencies until it's closed"
$ws.Range("B47").Value = ">>"
$ws.Range("B48").Value = ">>"

# Update view state to match the saved workbook: scrolled/selected position
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C49").Select()
